$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($ws, $row1, $row2, $colStart, $colEnd)
    $rng1 = $ws.Range($colStart + $row1 + ":" + $colEnd + $row1)
    $rng2 = $ws.Range($colStart + $row2 + ":" + $colEnd + $row2)
    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2
    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

Swap-RowRange $ws 136 137 "F" "V"
Swap-RowRange $ws 141 142 "F" "V"
Swap-RowRange $ws 147 148 "F" "V"
Swap-RowRange $ws 158 159 "F" "V"

# Append new row 161 (Burgos CF vs Alcorcon), copying formatting from row 160 first
$ws.Range("A160:V160").Copy($ws.Range("A161:V161"))

$ws.Range("A161").Value = 160
$ws.Range("B161").Value = "spain"
$ws.Range("C161").Value = "laliga2"
$ws.Range("D161").Value = "2023-2024"
$ws.Range("E161").Value = 45242.58333333334
$ws.Range("F161").Value = "Burgos CF"
$ws.Range("G161").Value = 4
$ws.Range("H161").Value = "Alcorcon"
$ws.Range("I161").Value = 2
$ws.Range("J161").Value = 2.08
$ws.Range("K161").Value = "05/11/2023 21:12"
$ws.Range("L161").Value = 2.03
$ws.Range("M161").Value = "12/11/2023 13:58"
$ws.Range("N161").Value = 3.21
$ws.Range("O161").Value = "05/11/2023 21:12"
$ws.Range("P161").Value = 3.18
$ws.Range("Q161").Value = "12/11/2023 13:58"
$ws.Range("R161").Value = 4.04
$ws.Range("S161").Value = "05/11/2023 21:12"
$ws.Range("T161").Value = 4.52
$ws.Range("U161").Value = "12/11/2023 13:58"
$ws.Range("V161").Value = "https://www.betexplorer.com/football/spain/laliga2/burgos-cf-alcorcon/04pxbtWt/"
